$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as exact text, avoiding Excel auto-converting
# numeric-looking strings (e.g. "1.000", "9.180") into numbers, and
# restore the default "Normal" style afterwards so no stray formatting
# is introduced relative to the original (unstyled) data cells.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Price column (D): forced as text to preserve exact formatting ---
Set-TextValue "D2" "30.263.48"
Set-TextValue "D3" "1.865.45"
Set-TextValue "D5" "234.98"
Set-TextValue "D7" "0.4678"
Set-TextValue "D8" "0.2835"
Set-TextValue "D9" "0.06507"
Set-TextValue "D10" "21.23"
Set-TextValue "D11" "0.07864"
Set-TextValue "D12" "97.25"
Set-TextValue "D13" "1.871.35"
Set-TextValue "D14" "5.092"
Set-TextValue "D15" "0.6725"
Set-TextValue "D16" "279.91"
Set-TextValue "D17" "30.257.12"
Set-TextValue "D18" "1.000"
Set-TextValue "D19" "5.483"
Set-TextValue "D20" "12.68"
Set-TextValue "D21" "2.123.90"
Set-TextValue "D25" "165.11"
Set-TextValue "D26" "9.181"
Set-TextValue "D27" "19.10"
Set-TextValue "D29" "1.377"
Set-TextValue "D30" "0.09636"
Set-TextValue "D31" "4.382"
Set-TextValue "D33" "4.093"
Set-TextValue "D34" "0.04699"
Set-TextValue "D35" "1.118"
Set-TextValue "D36" "0.7049"
Set-TextValue "D37" "2.727"
Set-TextValue "D38" "0.01850"
Set-TextValue "D39" "2.533"
Set-TextValue "D40" "6.246"
Set-TextValue "D41" "73.28"
Set-TextValue "D42" "1.942"
Set-TextValue "D43" "0.8459"
Set-TextValue "D44" "0.4167"
Set-TextValue "D46" "103.67"
Set-TextValue "D47" "7.166"
Set-TextValue "D48" "9.180"
Set-TextValue "D49" "936.93"
Set-TextValue "D50" "34.07"
Set-TextValue "D51" "0.1122"

# --- Coin name / link / volume columns (B, C, E): plain text assignment ---
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  -3.56%  "
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E40").Value = "  -5.82%  "
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("E49").Value = "  -5.00%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  -2.07%  "
